# CostsData.xlsx update:
#  - CO2 column header/units switched from lb/MMBtu to kg/MMBtu
#  - CO2 values for Coal, Gas-CC, Gas-CT converted from lb/MMBtu to kg/MMBtu
#  - Selected cell moved to J18

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header: CO2 unit changed from lb/MMBtu to kg/MMBtu
$ws.Range("H1").Value = "CO2 (kg/MMBtu)"

# Convert CO2 intensity values (lb/MMBtu -> kg/MMBtu)
$ws.Range("H9").Value = 91.761661600000011   # Coal
$ws.Range("H10").Value = 53.796011199999995  # Gas-CC
$ws.Range("H11").Value = 53.977448000000003  # Gas-CT

# Update the active selection
$ws.Range("J18").Select()
